# RSTK-9614-SYDATA-Qty Complete at Operation Reversal
# - Update B2 value (Lot Track project -> SYDATA1 lot track project)
# - Remove bold header formatting from A1:F1 (now plain/default style)
# - Widen column D
# - Update the active selection on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cell value that changed
$ws.Range("B2").Value = "Pro-SYDATA1 (Lot track)"

# Remove the special bold-font styling previously applied to the header
# row (A1:F1), returning those cells to the default/Normal style.
$ws.Range("A1:F1").Style = "Normal"

# Widen column D to fit the (unchanged) header text
$ws.Columns.Item(4).ColumnWidth = 7.6666666666666667

# Update sheet selection to match the saved view state (select the
# whole populated range)
$ws.Range("A1:K2").Select()
